# Update the public EPEX Spot prices workbook with the next day's data.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Prix Spot": add a new day column CD (03-sep) with hourly prices
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy formatting from the previous header cell (CC1) so the new header
# cell keeps the bold / bordered / centered style used by the others.
$wsPrix.Range("CC1").Copy($wsPrix.Range("CD1"))
$wsPrix.Range("CD1").Value = "03-sep"

$prixSpotValues = @{
    2  = 17.44
    3  = 15.65
    4  = 10.86
    5  = 8.46
    6  = 4.11
    7  = 8.94
    8  = 10.14
    9  = 21.34
    10 = 22.6
    11 = 16.37
    12 = 0
    13 = -0.01
    14 = -0.01
    15 = -0.02
    16 = -0.02
    17 = -0.01
    18 = -0.01
    19 = 0
    20 = 12.85
    21 = 56.79
    22 = 58.21
    23 = 53.41
    24 = 70.23
    25 = 56.33
}

foreach ($row in $prixSpotValues.Keys) {
    $wsPrix.Range("CD$row").Value = $prixSpotValues[$row]
}

# ---------------------------------------------------------------------
# Sheet "Gaz": append 2025-09-01 daily price
#
# The date-like text "2025-09-01" must be stored as plain text (as all
# the other date cells in column A are), not auto-converted by Excel
# into a date serial number. Entering it as a formula that returns a
# string, then pasting the result back as a value, keeps it as text
# without registering any new cell style.
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A79").Formula = '="2025-09-01"'
$wsGaz.Range("A79").Copy()
$wsGaz.Range("A79").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wsGaz.Range("B79").Value = 30.8

# ---------------------------------------------------------------------
# Sheet "CO2": append 2025-09-01 daily price
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")
$wsCO2.Range("A79").Formula = '="2025-09-01"'
$wsCO2.Range("A79").Copy()
$wsCO2.Range("A79").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$wsCO2.Range("B79").Value = 73.31
